$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'279.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.05%"
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'27.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.17%"
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'4.832"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.84%"
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'0.06412"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.02%"
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'7.046"
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'1.299"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'3.99%"
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.9046"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'2.41%"
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.1539"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.17%"
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.06243"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'23.30%"
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.07463"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.74%"
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.02923"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.14%"
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'0.08988"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.31%"
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'0.001583"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.94%"
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'0.0006436"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.22%"
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.006001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.15%"
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'3.484"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.78%"
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'3.306"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.32%"
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'2.233"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.72%"
$ws.Range("E19").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'0.1352"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.79%"
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'3.919"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.03%"
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'0.04401"
$ws.Range("D23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'0.1501"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'8.78%"
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'0.001176"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.07%"
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'0.004303"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'11.33%"
$ws.Range("E26").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'0.0001180"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'-1.73%"
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'0.0001655"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'-14.49%"
$ws.Range("E29").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.04083"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.23%"
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006649"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.43%"
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1406"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'19.51%"
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.002090"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-6.32%"
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.01106"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.66%"
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.00005542"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'6.44%"
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'1.628"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'9.82%"
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.01847"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-8.82%"
$ws.Range("E47").Style = "Normal"
